$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a Range to become its own standalone run (splitting it away
# from whatever run(s) it was sharing) by doing a throw-away text swap, then
# set it to the desired final text. This reliably produces a clean run
# boundary without having to hand-count placeholder string lengths.
# ---------------------------------------------------------------------------
function Split-And-SetText($rng, [string]$finalText) {
    $rng.Text = "X"
    $rng2 = $d.Range($rng.Start, $rng.End)
    $rng2.Text = $finalText
    return $d.Range($rng.Start, $rng.Start + $finalText.Length)
}

# ===========================================================================
# Edit 1: ". A cada uno de los elementos de decisión se le asigna un valor."
#   -> "... se le fija un valor." (with "fija" split into its own run, and
#      the document's _GoBack bookmark relocated to just after it)
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("se le asigna un valor") | Out-Null
$start = $rng.Start

$wordRange = $d.Range($start + 6, $start + 12)   # "asigna"
$wordRange.Font.Name = "Arial"
$wordRange.Text = "fija"
$afterFija = $wordRange.End

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmPoint = $d.Range($afterFija, $afterFija)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ===========================================================================
# Edit 2: "...calcula el peso usando la fórmula para asignar y le da..."
#   -> "...fórmula para fijar y le da..." (with "fijar" split into its own
#      run)
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("fórmula para asignar y le da") | Out-Null
$start = $rng.Start

$wordRange2 = $d.Range($start + 13, $start + 20)  # "asignar"
$wordRange2.Font.Name = "Arial"
$wordRange2.Text = "fijar"

# ===========================================================================
# Edit 3: "Libro de entrada para curso de postgrado, que recoge el nombre y
#   apellidos del estudiante, CI, curso de postgrado, tomo y folio del
#   Certificado de notas. "
#   -> keep text identical, but:
#     - merge the (previously two) italic runs "Libro de entrada para
#       curso " + "de postgrado," into a single italic run
#     - split what follows into two non-italic runs: "... tomo y folio
#       del " and "Certificado de notas. "
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("Libro de entrada para curso de postgrado, que recoge el nombre y apellidos del estudiante, CI, curso de postgrado, tomo y folio del Certificado de notas. ") | Out-Null
$start = $rng.Start

$italicText = "Libro de entrada para curso de postgrado,"
$italicRange = $d.Range($start, $start + $italicText.Length)
Split-And-SetText $italicRange $italicText | Out-Null
$afterItalic = $start + $italicText.Length

$part1Text = " que recoge el nombre y apellidos del estudiante, CI, curso de postgrado, tomo y folio del "
$part2Text = "Certificado de notas. "

$part2Range = $d.Range($afterItalic + $part1Text.Length, $afterItalic + $part1Text.Length + $part2Text.Length)
$part2Range.Font.Name = "Arial"
$part2Range.Text = $part2Text

Write-Output "Done"
